# Apply numeric profit-recalculation updates across all class sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7304.619
$ws.Range("I40").Value = 5850.2144
$ws.Range("J40").Value = 10213.429
$ws.Range("K40").Value = 5850.2144
$ws.Range("L40").Value = 10213.429
$ws.Range("M40").Value = -5675.2144
$ws.Range("N40").Value = -10563.429
$ws.Range("H112").Value = 2940.7273
$ws.Range("J112").Value = 2940.7273
$ws.Range("L112").Value = 8822.1819
$ws.Range("N112").Value = -11038.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3079
$ws.Range("I61").Value = 3131.6365
$ws.Range("K61").Value = 3131.6365
$ws.Range("M61").Value = -2919.6365
$ws.Range("H110").Value = 229038.36
$ws.Range("I110").Value = 314640.2
$ws.Range("K110").Value = 314640.2
$ws.Range("M110").Value = -312595.2
$ws.Range("H122").Value = 3980.5715
$ws.Range("I122").Value = 2260.923
$ws.Range("K122").Value = 6782.768999999999
$ws.Range("M122").Value = -4332.768999999999
$ws.Range("H136").Value = 3079
$ws.Range("I136").Value = 3131.6365
$ws.Range("K136").Value = 9394.9095
$ws.Range("M136").Value = -6844.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H86").Value = 775270.25
$ws.Range("I86").Value = 1310076.6
$ws.Range("K86").Value = 1310076.6
$ws.Range("M86").Value = -1308953.6
$ws.Range("H89").Value = 775270.25
$ws.Range("I89").Value = 1310076.6
$ws.Range("K89").Value = 6550383
$ws.Range("M89").Value = -6544767
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3354
$ws.Range("I16").Value = 3443.5
$ws.Range("J16").Value = 2906.5
$ws.Range("K16").Value = 3443.5
$ws.Range("L16").Value = 2906.5
$ws.Range("M16").Value = -3156.5
$ws.Range("N16").Value = -3480.5
$ws.Range("H58").Value = 1695.5217
$ws.Range("I58").Value = 1437.8889
$ws.Range("K58").Value = 1437.8889
$ws.Range("M58").Value = -1234.8889
$ws.Range("H59").Value = 30566.182
$ws.Range("I59").Value = 30000
$ws.Range("J59").Value = 30622.8
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 30622.8
$ws.Range("M59").Value = -28855
$ws.Range("N59").Value = -32912.8
$ws.Range("H60").Value = 37142.855
$ws.Range("J60").Value = 48000
$ws.Range("L60").Value = 48000
$ws.Range("N60").Value = -49022
$ws.Range("H107").Value = 671.04
$ws.Range("I107").Value = 544.5
$ws.Range("K107").Value = 544.5
$ws.Range("M107").Value = 1375.5
$ws.Range("H113").Value = 3354
$ws.Range("I113").Value = 3443.5
$ws.Range("J113").Value = 2906.5
$ws.Range("K113").Value = 3443.5
$ws.Range("L113").Value = 2906.5
$ws.Range("M113").Value = -1273.5
$ws.Range("N113").Value = -7246.5
$ws.Range("H125").Value = 98000
$ws.Range("J125").Value = 98000
$ws.Range("L125").Value = 98000
$ws.Range("N125").Value = -102920
$ws.Range("H132").Value = 1998.3572
$ws.Range("I132").Value = 1620.9524
$ws.Range("K132").Value = 4862.857199999999
$ws.Range("M132").Value = -2332.857199999999
$ws.Range("H136").Value = 1695.5217
$ws.Range("I136").Value = 1437.8889
$ws.Range("K136").Value = 4313.6667
$ws.Range("M136").Value = -1763.6667
$ws.Range("H141").Value = 315056.75
$ws.Range("J141").Value = 327229.3
$ws.Range("L141").Value = 327229.3
$ws.Range("N141").Value = -337589.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360
$ws.Range("H138").Value = 3309.3333
$ws.Range("I138").Value = 2714.25
$ws.Range("J138").Value = 4499.5
$ws.Range("K138").Value = 8142.75
$ws.Range("L138").Value = 13498.5
$ws.Range("M138").Value = -3002.75
$ws.Range("N138").Value = -23778.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 912840.3
$ws.Range("I80").Value = 718603.6
$ws.Range("J80").Value = 1252754.5
$ws.Range("K80").Value = 718603.6
$ws.Range("L80").Value = 1252754.5
$ws.Range("M80").Value = -717605.6
$ws.Range("N80").Value = -1254750.5
$ws.Range("H83").Value = 912840.3
$ws.Range("I83").Value = 718603.6
$ws.Range("J83").Value = 1252754.5
$ws.Range("K83").Value = 3593018
$ws.Range("L83").Value = 6263772.5
$ws.Range("M83").Value = -3588026
$ws.Range("N83").Value = -6273756.5
$ws.Range("H97").Value = 663.5833
$ws.Range("I97").Value = 717.64703
$ws.Range("K97").Value = 717.64703
$ws.Range("M97").Value = -221.64703
$ws.Range("H102").Value = 3956.3684
$ws.Range("I102").Value = 1988.1666
$ws.Range("K102").Value = 1988.1666
$ws.Range("M102").Value = -366.1666
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H113").Value = 386995
$ws.Range("I113").Value = 589618.8
$ws.Range("K113").Value = 589618.8
$ws.Range("M113").Value = -587448.8
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -79900
$ws.Range("H124").Value = 85800
$ws.Range("J124").Value = 85800
$ws.Range("L124").Value = 85800
$ws.Range("N124").Value = -95620
$ws.Range("H132").Value = 36687.332
$ws.Range("I132").Value = 3607.4814
$ws.Range("J132").Value = 334406
$ws.Range("K132").Value = 10822.4442
$ws.Range("L132").Value = 1003218
$ws.Range("M132").Value = -8292.4442
$ws.Range("N132").Value = -1008278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5005.7085
$ws.Range("I40").Value = 4217.737
$ws.Range("K40").Value = 4217.737
$ws.Range("M40").Value = -4081.737
$ws.Range("H61").Value = 2853.7144
$ws.Range("I61").Value = 2946.55
$ws.Range("K61").Value = 2946.55
$ws.Range("M61").Value = -2744.55
$ws.Range("H113").Value = 2853.7144
$ws.Range("I113").Value = 2946.55
$ws.Range("K113").Value = 2946.55
$ws.Range("M113").Value = -776.5500000000002
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35718636
$ws.Range("I122").Value = 47623428
$ws.Range("J122").Value = 4271.143
$ws.Range("K122").Value = 142870284
$ws.Range("L122").Value = 12813.429
$ws.Range("M122").Value = -142867834
$ws.Range("N122").Value = -17713.429
$ws.Range("H136").Value = 79515.766
$ws.Range("I136").Value = 2720.35
$ws.Range("K136").Value = 8161.049999999999
$ws.Range("M136").Value = -5611.049999999999
